$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 0
$ws.Range("K6").Value2 = 0
$ws.Range("M6").Value2 = $null
$ws.Range("H8").Value2 = 258.2
$ws.Range("I8").Value2 = 181.33333
$ws.Range("K8").Value2 = 543.99999
$ws.Range("M8").Value2 = -404.99999
$ws.Range("H13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("N13").Value2 = $null
$ws.Range("H76").Value2 = 4124.8335
$ws.Range("I76").Value2 = 2500
$ws.Range("J76").Value2 = 4937.25
$ws.Range("K76").Value2 = 2500
$ws.Range("L76").Value2 = 4937.25
$ws.Range("M76").Value2 = -2185
$ws.Range("N76").Value2 = -5567.25
$ws.Range("H79").Value2 = 4124.8335
$ws.Range("I79").Value2 = 2500
$ws.Range("J79").Value2 = 4937.25
$ws.Range("K79").Value2 = 2500
$ws.Range("L79").Value2 = 4937.25
$ws.Range("M79").Value2 = -1408
$ws.Range("N79").Value2 = -7121.25
$ws.Range("H125").Value2 = 2658.7144
$ws.Range("I125").Value2 = 5781
$ws.Range("J125").Value2 = 2138.3333
$ws.Range("K125").Value2 = 52029
$ws.Range("L125").Value2 = 19244.9997
$ws.Range("M125").Value2 = -49569
$ws.Range("N125").Value2 = -24164.9997
$ws.Range("H127").Value2 = 16951314
$ws.Range("I127").Value2 = 649.8
$ws.Range("J127").Value2 = 18520818
$ws.Range("K127").Value2 = 1949.4
$ws.Range("L127").Value2 = 55562454
$ws.Range("M127").Value2 = 3010.6
$ws.Range("N127").Value2 = -55572374
$ws.Range("H129").Value2 = 947.76
$ws.Range("I129").Value2 = 606.5
$ws.Range("J129").Value2 = 1012.7619
$ws.Range("K129").Value2 = 1819.5
$ws.Range("L129").Value2 = 3038.2857
$ws.Range("M129").Value2 = 3180.5
$ws.Range("N129").Value2 = -13038.2857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value2 = 8000
$ws.Range("I3").Value2 = 0
$ws.Range("K3").Value2 = 0
$ws.Range("M3").Value2 = $null
$ws.Range("H17").Value2 = 0
$ws.Range("J17").Value2 = 0
$ws.Range("L17").Value2 = 0
$ws.Range("N17").Value2 = $null
$ws.Range("H19").Value2 = 500
$ws.Range("I19").Value2 = 500
$ws.Range("K19").Value2 = 500
$ws.Range("M19").Value2 = -271
$ws.Range("H22").Value2 = 517.2
$ws.Range("I22").Value2 = 517.2
$ws.Range("K22").Value2 = 517.2
$ws.Range("M22").Value2 = -218.2
$ws.Range("H61").Value2 = 1961.1428
$ws.Range("I61").Value2 = 1348.5333
$ws.Range("J61").Value2 = 2928.4211
$ws.Range("K61").Value2 = 1348.5333
$ws.Range("L61").Value2 = 2928.4211
$ws.Range("M61").Value2 = -1136.5333
$ws.Range("N61").Value2 = -3352.4211
$ws.Range("H102").Value2 = 45797.348
$ws.Range("I102").Value2 = 78736.16
$ws.Range("J102").Value2 = 2976.9
$ws.Range("K102").Value2 = 78736.16
$ws.Range("L102").Value2 = 2976.9
$ws.Range("M102").Value2 = -77114.16
$ws.Range("N102").Value2 = -6220.9
$ws.Range("H136").Value2 = 1961.1428
$ws.Range("I136").Value2 = 1348.5333
$ws.Range("J136").Value2 = 2928.4211
$ws.Range("K136").Value2 = 4045.5999
$ws.Range("L136").Value2 = 8785.263300000001
$ws.Range("M136").Value2 = -1495.5999
$ws.Range("N136").Value2 = -13885.2633

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value2 = 0
$ws.Range("J48").Value2 = 0
$ws.Range("L48").Value2 = 0
$ws.Range("N48").Value2 = $null
$ws.Range("H55").Value2 = 10034.125
$ws.Range("I55").Value2 = 3868.25
$ws.Range("J55").Value2 = 16200
$ws.Range("K55").Value2 = 3868.25
$ws.Range("L55").Value2 = 16200
$ws.Range("M55").Value2 = -3553.25
$ws.Range("N55").Value2 = -16830
$ws.Range("H99").Value2 = 11720.77
$ws.Range("I99").Value2 = 4441.4546
$ws.Range("J99").Value2 = 51757
$ws.Range("K99").Value2 = 4441.4546
$ws.Range("L99").Value2 = 51757
$ws.Range("M99").Value2 = -2943.4546
$ws.Range("N99").Value2 = -54753
$ws.Range("H126").Value2 = 11720.77
$ws.Range("I126").Value2 = 4441.4546
$ws.Range("J126").Value2 = 51757
$ws.Range("K126").Value2 = 13324.3638
$ws.Range("L126").Value2 = 155271
$ws.Range("M126").Value2 = -10854.3638
$ws.Range("N126").Value2 = -160211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value2 = 144.44444
$ws.Range("K4").Value2 = 433.33332
$ws.Range("M4").Value2 = -321.33332
$ws.Range("H6").Value2 = 79.8
$ws.Range("I6").Value2 = 79.8
$ws.Range("J6").Value2 = 0
$ws.Range("K6").Value2 = 239.4
$ws.Range("L6").Value2 = 0
$ws.Range("M6").Value2 = -126.4
$ws.Range("N6").Value2 = $null
$ws.Range("H9").Value2 = 1000
$ws.Range("J9").Value2 = 1000
$ws.Range("L9").Value2 = 3000
$ws.Range("N9").Value2 = -3448
$ws.Range("H10").Value2 = 55.833332
$ws.Range("J10").Value2 = 0
$ws.Range("L10").Value2 = 0
$ws.Range("N10").Value2 = $null
$ws.Range("H13").Value2 = 100
$ws.Range("I13").Value2 = 100
$ws.Range("K13").Value2 = 300
$ws.Range("M13").Value2 = -132
$ws.Range("H54").Value2 = 2500
$ws.Range("J54").Value2 = 2500
$ws.Range("L54").Value2 = 7500
$ws.Range("N54").Value2 = -8618
$ws.Range("H87").Value2 = 8645.454
$ws.Range("I87").Value2 = 5050
$ws.Range("J87").Value2 = 12960
$ws.Range("K87").Value2 = 15150
$ws.Range("L87").Value2 = 38880
$ws.Range("M87").Value2 = -13902
$ws.Range("N87").Value2 = -41376
$ws.Range("H90").Value2 = 8645.454
$ws.Range("I90").Value2 = 5050
$ws.Range("J90").Value2 = 12960
$ws.Range("K90").Value2 = 45450
$ws.Range("L90").Value2 = 116640
$ws.Range("M90").Value2 = -39210
$ws.Range("N90").Value2 = -129120
$ws.Range("H114").Value2 = 754.25
$ws.Range("I114").Value2 = 228.33333
$ws.Range("K114").Value2 = 684.99999
$ws.Range("M114").Value2 = 2569.00001
$ws.Range("H122").Value2 = 14925.429
$ws.Range("J122").Value2 = 20759.6
$ws.Range("L122").Value2 = 186836.4
$ws.Range("N122").Value2 = -191736.4
$ws.Range("H134").Value2 = 2585.7576
$ws.Range("I134").Value2 = 2048.889
$ws.Range("J134").Value2 = 2787.0833
$ws.Range("K134").Value2 = 6146.667
$ws.Range("L134").Value2 = 8361.249899999999
$ws.Range("M134").Value2 = -1076.667
$ws.Range("N134").Value2 = -18501.2499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 5000000
$ws.Range("J5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("N5").Value2 = $null
$ws.Range("H9").Value2 = 747
$ws.Range("I9").Value2 = 747
$ws.Range("K9").Value2 = 747
$ws.Range("M9").Value2 = -577
$ws.Range("H19").Value2 = 46002
$ws.Range("I19").Value2 = 9999
$ws.Range("J19").Value2 = 70004
$ws.Range("K19").Value2 = 9999
$ws.Range("L19").Value2 = 70004
$ws.Range("M19").Value2 = -9711
$ws.Range("N19").Value2 = -70580

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value2 = 282437.5
$ws.Range("I2").Value2 = 300600
$ws.Range("J2").Value2 = 10000
$ws.Range("K2").Value2 = 300600
$ws.Range("L2").Value2 = 10000
$ws.Range("M2").Value2 = -300488
$ws.Range("N2").Value2 = -10224
$ws.Range("H13").Value2 = 46674.668
$ws.Range("I13").Value2 = 10
$ws.Range("K13").Value2 = 10
$ws.Range("M13").Value2 = 130

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value2 = 1013379.7
$ws.Range("I2").Value2 = 1120444.1
$ws.Range("K2").Value2 = 1120444.1
$ws.Range("M2").Value2 = -1120332.1
$ws.Range("H13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("L13").Value2 = 0
$ws.Range("N13").Value2 = $null
$ws.Range("H22").Value2 = 5906.5
$ws.Range("J22").Value2 = 9800
$ws.Range("L22").Value2 = 9800
$ws.Range("N22").Value2 = -10386
$ws.Range("H23").Value2 = 951.25
$ws.Range("I23").Value2 = 351.66666
$ws.Range("J23").Value2 = 2750
$ws.Range("K23").Value2 = 351.66666
$ws.Range("L23").Value2 = 2750
$ws.Range("M23").Value2 = -122.66666
$ws.Range("N23").Value2 = -3208
$ws.Range("H95").Value2 = 0
$ws.Range("J95").Value2 = 0
$ws.Range("L95").Value2 = 0
$ws.Range("N95").Value2 = $null
$ws.Range("H132").Value2 = 3719.9
$ws.Range("I132").Value2 = 4251.385
$ws.Range("K132").Value2 = 12754.155
$ws.Range("M132").Value2 = -10224.155
